# 16-Apr-2024: GUI implemented to load the configuration file and specify
# how many test paper to be generated.
#
# The config sheet gains two new option rows describing how the "test
# paper" and "marksheet" file names are produced (the extension is now
# appended by the program instead of being baked into the stored value),
# plus explanatory notes in columns C/D for the online/paper-test rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8 (candidates data) : new note in column C ---------------------
$ws.Range("C8").Value = "; online test"

# --- Row 9 (test results) : new note in column C ------------------------
$ws.Range("C9").Value = "; online test"

# --- Row 10 (language) : new note in column C ----------------------------
$ws.Range("C10").Value = "; paper test"

# --- Row 11 (test paper) : value trimmed to base name, notes added ------
$ws.Range("B11").Value = "testpaper"
$ws.Range("C11").Value = "; paper test"
$ws.Range("D11").Value = "the file type (pdf) will be added by the program"

# --- Row 12 (marksheet) : value trimmed to base name, notes added -------
$ws.Range("B12").Value = "marksheet"
$ws.Range("C12").Value = "; paper test"
$ws.Range("D12").Value = "the file type (xlsx) will be added by the program"

# Reset the active selection back to the top-left cell (the saved file no
# longer pins the selection on B10).
$ws.Range("A1").Select()
